$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the "level" column (B) - not part of the rebuilt storeItemDefine schema
$ws.Columns.Item(2).Delete()

# Drop the "iconResource"/"desc" columns (originally F:G, now E:F after the delete above)
$ws.Range("E1:F5").Delete()

# Re-index the ID column to be 0-based (was 24-27)
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3

# New width for the maxPrice column (D)
$ws.Columns.Item(4).ColumnWidth = 11.428571428571429

# Match the author's final selection
$ws.Range("E5").Select() | Out-Null
